# Regenerate merged AHB files
# - rename the header row's "_old"/"_new" suffixed labels to "_FV2404"/"_FV2410"
# - (re)create the Table1 list-object over the full data range
# - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header labels (row 1) ---------------------------------
# A1:J1  "<Name>_old" -> "<Name>_FV2404"
# K1     "diff"       -> unchanged
# L1:U1  "<Name>_new" -> "<Name>_FV2410"
$headersFv2404 = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)
$headersFv2410 = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

# columns A..J -> _old -> _FV2404
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headersFv2404[$i]
}
# column K ("diff") is unchanged

# columns L..U -> _new -> _FV2410
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $headersFv2410[$i]
}

# --- 2. Turn the data range into an Excel Table ("Table1") -----------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U60"), [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row -----------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true | Out-Null
